$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("first_eval")

# Row 2
$ws.Range("B2").Value = 0.2147957804815247
$ws.Range("C2").Value = 0.3580819848962541
$ws.Range("D2").Value = 0.2088379301833887
$ws.Range("E2").Value = 0.4569878884427777
$ws.Range("F2").Value = 0.4185882721160311
$ws.Range("G2").Value = 14

# Row 3
$ws.Range("B3").Value = 0.245473867636748
$ws.Range("C3").Value = 0.3683843928281823
$ws.Range("D3").Value = 0.2358514199634523
$ws.Range("E3").Value = 0.4856453644002507
$ws.Range("F3").Value = 0.4361500127557235
$ws.Range("G3").Value = 13

# Row 4
$ws.Range("B4").Value = 0.301463903713677
$ws.Range("C4").Value = 0.3567050099425655
$ws.Range("D4").Value = 0.2583510523008795
$ws.Range("E4").Value = 0.5082824532687308
$ws.Range("F4").Value = 0.4274285484895888
$ws.Range("G4").Value = 12

# Row 5
$ws.Range("B5").Value = 0.2900334324511264
$ws.Range("C5").Value = 0.3439303588265401
$ws.Range("D5").Value = 0.1874043702602189
$ws.Range("E5").Value = 0.4329022640968963
$ws.Range("F5").Value = 0.3370659818980856
$ws.Range("G5").Value = 11

# Row 6
$ws.Range("B6").Value = 0.2833661826906564
$ws.Range("C6").Value = 0.2909807410604855
$ws.Range("D6").Value = 0.2099604875429734
$ws.Range("E6").Value = 0.4582144558424291
$ws.Range("F6").Value = 0.3795671424286398
$ws.Range("G6").Value = 10

# Row 7
$ws.Range("B7").Value = 0.3070978867771534
$ws.Range("C7").Value = 0.3733636115588901
$ws.Range("D7").Value = 0.2272551859511885
$ws.Range("E7").Value = 0.476712896774556
$ws.Range("F7").Value = 0.3867354821117653
$ws.Range("G7").Value = 9

# Row 8
$ws.Range("B8").Value = 0.3125739411747067
$ws.Range("C8").Value = 0.4101416856701894
$ws.Range("D8").Value = 0.3161064439746084
$ws.Range("E8").Value = 0.5622334425971194
$ws.Range("F8").Value = 0.499604386659364
$ws.Range("G8").Value = 8

# Row 9
$ws.Range("B9").Value = 0.2627214152773011
$ws.Range("C9").Value = 0.3235434349789264
$ws.Range("D9").Value = 0.1669924024756296
$ws.Range("E9").Value = 0.4086470389904099
$ws.Range("F9").Value = 0.3380801243818421
$ws.Range("G9").Value = 7

# Row 10
$ws.Range("B10").Value = 0.3023013051968123
$ws.Range("C10").Value = 0.3856149542229345
$ws.Range("D10").Value = 0.3827274330185637
$ws.Range("E10").Value = 0.618649685216572
$ws.Range("F10").Value = 0.5912779588939882
$ws.Range("G10").Value = 6

# Row 11
$ws.Range("B11").Value = 0.302435154479656
$ws.Range("C11").Value = 0.4668197381621303
$ws.Range("D11").Value = 0.3675323742590614
$ws.Range("E11").Value = 0.6062444838999044
$ws.Range("F11").Value = 0.5874365408215683
$ws.Range("G11").Value = 5
